$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D7"  = -7.632
    "C8"  = -12.679
    "C10" = -12.937
    "C12" = -11.207
    "D15" = -8.378000000000002
    "C18" = -14.011
    "D18" = -8.456
    "D20" = -7.517000000000001
    "D29" = -7.292
    "D30" = -7.141
    "D31" = -7.938
    "C37" = -13.243
    "D40" = -7.641
    "D50" = -8.008000000000001
    "C55" = -13.958
    "C68" = -11.167
    "D68" = -6.879
    "D76" = -7.312
    "C77" = -13.117
    "C78" = -13.214
    "C81" = -13.094
    "C82" = -11.737
    "D87" = -8.339
    "D88" = -8.177000000000001
    "D96" = -7.267
    "D98" = -8.242000000000001
    "D101" = -7.885999999999998
    "D102" = -8.036
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
